# Re-applies the refreshed "previous copy of ful-path.csv" export onto
# the quadratic-svm-score sheet:
#   - column A narrows to fit the short numeric "Row" labels
#   - header row (A1:C1) keeps its text formatting
#   - row 2's "Row" label now points at the top-scoring genome
#     (even_MAG-GUT68245.fa, shared-string index 50) instead of the
#     first alphabetical one, and its score is refreshed from 1 to the
#     real predicted value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column A (was sized for the long "even_MAG-GUT....fa" labels,
# now only needs to fit small row numbers / short labels).
$ws.Columns("A:A").ColumnWidth = 4.166666666666667

# Header row - re-apply text format so it keeps its shared style.
$ws.Range("A1:C1").NumberFormat = "@"

# Row 2 data refresh.
$ws.Range("A2").Value2 = "even_MAG-GUT68245.fa"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").Value2 = 13.031067311341907
$ws.Range("C2").Value2 = 1
